$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 296.875
$ws.Range("I5").Value = 145
$ws.Range("K5").Value = 145
$ws.Range("M5").Value = -30
$ws.Range("H42").Value = 795.2778
$ws.Range("I42").Value = 289
$ws.Range("J42").Value = 1807.8334
$ws.Range("K42").Value = 867
$ws.Range("L42").Value = 5423.5002
$ws.Range("M42").Value = -637
$ws.Range("N42").Value = -5883.5002
$ws.Range("H51").Value = 2740.75
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2740.75
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2740.75
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -3708.75
$ws.Range("H69").Value = 8996.723
$ws.Range("I69").Value = 4749.5
$ws.Range("J69").Value = 9527.625
$ws.Range("K69").Value = 14248.5
$ws.Range("L69").Value = 28582.875
$ws.Range("M69").Value = -13374.5
$ws.Range("N69").Value = -30330.875
$ws.Range("H70").Value = 1870.4445
$ws.Range("I70").Value = 950
$ws.Range("K70").Value = 2850
$ws.Range("M70").Value = -2580
$ws.Range("H72").Value = 8996.723
$ws.Range("I72").Value = 4749.5
$ws.Range("J72").Value = 9527.625
$ws.Range("K72").Value = 42745.5
$ws.Range("L72").Value = 85748.625
$ws.Range("M72").Value = -38377.5
$ws.Range("N72").Value = -94484.625
$ws.Range("H73").Value = 1870.4445
$ws.Range("I73").Value = 950
$ws.Range("K73").Value = 2850
$ws.Range("M73").Value = -1914
$ws.Range("H82").Value = 3262.5
$ws.Range("I82").Value = 3683.3333
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 11049.9999
$ws.Range("L82").Value = 6000
$ws.Range("M82").Value = -10643.9999
$ws.Range("N82").Value = -6812
$ws.Range("H85").Value = 3262.5
$ws.Range("I85").Value = 3683.3333
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 11049.9999
$ws.Range("L85").Value = 6000
$ws.Range("M85").Value = -9645.999899999999
$ws.Range("N85").Value = -8808
$ws.Range("H112").Value = 6095
$ws.Range("J112").Value = 6174.537
$ws.Range("L112").Value = 18523.611
$ws.Range("N112").Value = -20739.611
$ws.Range("H137").Value = 957811.2
$ws.Range("I137").Value = 2807.6316
$ws.Range("J137").Value = 2167482.2
$ws.Range("K137").Value = 8422.8948
$ws.Range("L137").Value = 6502446.600000001
$ws.Range("M137").Value = -5872.8948
$ws.Range("N137").Value = -6507546.600000001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3334033
$ws.Range("I61").Value = 3334033
$ws.Range("K61").Value = 3334033
$ws.Range("M61").Value = -3333821
$ws.Range("H74").Value = 4337.5293
$ws.Range("I74").Value = 1255.1428
$ws.Range("J74").Value = 6495.2
$ws.Range("K74").Value = 1255.1428
$ws.Range("L74").Value = 6495.2
$ws.Range("M74").Value = -381.1428000000001
$ws.Range("N74").Value = -8243.200000000001
$ws.Range("H77").Value = 4337.5293
$ws.Range("I77").Value = 1255.1428
$ws.Range("J77").Value = 6495.2
$ws.Range("K77").Value = 6275.714
$ws.Range("L77").Value = 32476
$ws.Range("M77").Value = -1907.714
$ws.Range("N77").Value = -41212
$ws.Range("H102").Value = 37149.5
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H110").Value = 1940.138
$ws.Range("I110").Value = 1565.5834
$ws.Range("J110").Value = 3738
$ws.Range("K110").Value = 1565.5834
$ws.Range("L110").Value = 3738
$ws.Range("M110").Value = 479.4166
$ws.Range("N110").Value = -7828
$ws.Range("H135").Value = 89986.336
$ws.Range("J135").Value = 89986.336
$ws.Range("L135").Value = 89986.336
$ws.Range("N135").Value = -100126.336
$ws.Range("H136").Value = 3334033
$ws.Range("I136").Value = 3334033
$ws.Range("K136").Value = 10002099
$ws.Range("M136").Value = -9999549

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6230.4546
$ws.Range("I105").Value = 6153.5
$ws.Range("J105").Value = 7000
$ws.Range("K105").Value = 6153.5
$ws.Range("L105").Value = 7000
$ws.Range("M105").Value = -4406.5
$ws.Range("N105").Value = -10494
$ws.Range("H134").Value = 694791.9399999999
$ws.Range("I134").Value = 746381.75
$ws.Range("J134").Value = 511361.44
$ws.Range("K134").Value = 2239145.25
$ws.Range("L134").Value = 1534084.32
$ws.Range("M134").Value = -2236610.25
$ws.Range("N134").Value = -1539154.32

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11052.631
$ws.Range("I31").Value = 3932.1562
$ws.Range("J31").Value = 27328
$ws.Range("K31").Value = 3932.1562
$ws.Range("L31").Value = 27328
$ws.Range("M31").Value = -3637.1562
$ws.Range("N31").Value = -27918
$ws.Range("H34").Value = 11052.631
$ws.Range("I34").Value = 3932.1562
$ws.Range("J34").Value = 27328
$ws.Range("K34").Value = 3932.1562
$ws.Range("L34").Value = 27328
$ws.Range("M34").Value = -3730.1562
$ws.Range("N34").Value = -27732
$ws.Range("H87").Value = 123553.336
$ws.Range("J87").Value = 123553.336
$ws.Range("L87").Value = 123553.336
$ws.Range("N87").Value = -125925.336
$ws.Range("H90").Value = 123553.336
$ws.Range("J90").Value = 123553.336
$ws.Range("L90").Value = 370660.008
$ws.Range("N90").Value = -382516.008
$ws.Range("H107").Value = 973.58826
$ws.Range("I107").Value = 870.13336
$ws.Range("J107").Value = 1749.5
$ws.Range("K107").Value = 870.13336
$ws.Range("L107").Value = 1749.5
$ws.Range("M107").Value = 1049.86664
$ws.Range("N107").Value = -5589.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H140").Value = 2613.4443
$ws.Range("I140").Value = 1878.5
$ws.Range("J140").Value = 4083.3333
$ws.Range("K140").Value = 5635.5
$ws.Range("L140").Value = 12249.9999
$ws.Range("M140").Value = -455.5
$ws.Range("N140").Value = -22609.9999

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3500.0833
$ws.Range("I68").Value = 3713
$ws.Range("J68").Value = 3202
$ws.Range("K68").Value = 3713
$ws.Range("L68").Value = 3202
$ws.Range("M68").Value = -2964
$ws.Range("N68").Value = -4700
$ws.Range("H71").Value = 3500.0833
$ws.Range("I71").Value = 3713
$ws.Range("J71").Value = 3202
$ws.Range("K71").Value = 18565
$ws.Range("L71").Value = 16010
$ws.Range("M71").Value = -14821
$ws.Range("N71").Value = -23498
$ws.Range("H99").Value = 33939
$ws.Range("I99").Value = 33939
$ws.Range("K99").Value = 33939
$ws.Range("M99").Value = -30944
$ws.Range("H136").Value = 3992
$ws.Range("I136").Value = 2700
$ws.Range("J136").Value = 5499.3335
$ws.Range("K136").Value = 8100
$ws.Range("L136").Value = 16498.0005
$ws.Range("M136").Value = -5550
$ws.Range("N136").Value = -21598.0005

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 43342
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 9150083
$ws.Range("I132").Value = 10063691
$ws.Range("K132").Value = 30191073
$ws.Range("M132").Value = -30188543
$ws.Range("H136").Value = 13726.675
$ws.Range("I136").Value = 13764.789
$ws.Range("K136").Value = 41294.367
$ws.Range("M136").Value = -38744.367
